$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 24, shifting the existing rows 24:56 down to 25:57.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new weekly record. All the
# "descriptive" columns are identical to the (now shifted-down) row 25
# record, except the date, volume and origin which change for this
# new entry.
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44469
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108002
$ws.Range("J24").Value = "Mango"
$ws.Range("K24").Value = "Sin especificar"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 8500
$ws.Range("O24").Value = 9000
$ws.Range("P24").Value = 8750
$ws.Range("Q24").Value = "$/bandeja 4 kilos"
$ws.Range("R24").Value = "Brasil"
$ws.Range("S24").Value = 2188
$ws.Range("T24").Value = 4
